$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("ALC").Range("H17").Value = 475.85
$wb.Worksheets.Item("ALC").Range("J17").Value = 475.85
$wb.Worksheets.Item("ALC").Range("L17").Value = 1427.55
$wb.Worksheets.Item("ALC").Range("N17").Value = -1763.55

$wb.Worksheets.Item("ALC").Range("H33").Value = 429.05554
$wb.Worksheets.Item("ALC").Range("I33").Value = 389.58823
$wb.Worksheets.Item("ALC").Range("K33").Value = 389.58823
$wb.Worksheets.Item("ALC").Range("M33").Value = -160.58823

$wb.Worksheets.Item("ALC").Range("H86").Value = 3899.75
$wb.Worksheets.Item("ALC").Range("I86").Value = 5649.5
$wb.Worksheets.Item("ALC").Range("J86").Value = 2150
$wb.Worksheets.Item("ALC").Range("K86").Value = 5649.5
$wb.Worksheets.Item("ALC").Range("L86").Value = 2150
$wb.Worksheets.Item("ALC").Range("M86").Value = -4526.5
$wb.Worksheets.Item("ALC").Range("N86").Value = -4396

$wb.Worksheets.Item("ALC").Range("H89").Value = 3899.75
$wb.Worksheets.Item("ALC").Range("I89").Value = 5649.5
$wb.Worksheets.Item("ALC").Range("J89").Value = 2150
$wb.Worksheets.Item("ALC").Range("K89").Value = 28247.5
$wb.Worksheets.Item("ALC").Range("L89").Value = 10750
$wb.Worksheets.Item("ALC").Range("M89").Value = -22631.5
$wb.Worksheets.Item("ALC").Range("N89").Value = -21982

$wb.Worksheets.Item("ALC").Range("H113").Value = 12006.546
$wb.Worksheets.Item("ALC").Range("I113").Value = 6799.4
$wb.Worksheets.Item("ALC").Range("K113").Value = 6799.4
$wb.Worksheets.Item("ALC").Range("M113").Value = -3545.4

$wb.Worksheets.Item("ALC").Range("H137").Value = 14921.454
$wb.Worksheets.Item("ALC").Range("I137").Value = 7154
$wb.Worksheets.Item("ALC").Range("K137").Value = 21462
$wb.Worksheets.Item("ALC").Range("M137").Value = -18912

$wb.Worksheets.Item("ALC").Range("H138").Value = 4577.3125
$wb.Worksheets.Item("ALC").Range("I138").Value = 1846.75
$wb.Worksheets.Item("ALC").Range("K138").Value = 5540.25
$wb.Worksheets.Item("ALC").Range("M138").Value = -400.25

$wb.Worksheets.Item("ALC").Range("H140").Value = 0
$wb.Worksheets.Item("ALC").Range("J140").Value = 0
$wb.Worksheets.Item("ALC").Range("L140").Value = 0
$wb.Worksheets.Item("ALC").Range("N140").Value = $null

$wb.Worksheets.Item("ARM").Range("H5").Value = 1375
$wb.Worksheets.Item("ARM").Range("I5").Value = 250
$wb.Worksheets.Item("ARM").Range("K5").Value = 250
$wb.Worksheets.Item("ARM").Range("M5").Value = -138

$wb.Worksheets.Item("ARM").Range("H32").Value = 703.2917
$wb.Worksheets.Item("ARM").Range("I32").Value = 471.73914
$wb.Worksheets.Item("ARM").Range("K32").Value = 471.73914
$wb.Worksheets.Item("ARM").Range("M32").Value = -184.73914

$wb.Worksheets.Item("ARM").Range("H33").Value = 6254.091
$wb.Worksheets.Item("ARM").Range("I33").Value = 6254.091
$wb.Worksheets.Item("ARM").Range("J33").Value = 0
$wb.Worksheets.Item("ARM").Range("K33").Value = 6254.091
$wb.Worksheets.Item("ARM").Range("L33").Value = 0
$wb.Worksheets.Item("ARM").Range("M33").Value = -5925.091
$wb.Worksheets.Item("ARM").Range("N33").Value = $null

$wb.Worksheets.Item("ARM").Range("H110").Value = 4345.696
$wb.Worksheets.Item("ARM").Range("I110").Value = 2252.1538
$wb.Worksheets.Item("ARM").Range("K110").Value = 2252.1538
$wb.Worksheets.Item("ARM").Range("M110").Value = -207.1538

$wb.Worksheets.Item("ARM").Range("H122").Value = 5450.25
$wb.Worksheets.Item("ARM").Range("I122").Value = 4022.889
$wb.Worksheets.Item("ARM").Range("K122").Value = 12068.667
$wb.Worksheets.Item("ARM").Range("M122").Value = -9618.667000000001

$wb.Worksheets.Item("BSM").Range("H4").Value = 1375
$wb.Worksheets.Item("BSM").Range("I4").Value = 250
$wb.Worksheets.Item("BSM").Range("K4").Value = 250
$wb.Worksheets.Item("BSM").Range("M4").Value = -135

$wb.Worksheets.Item("BSM").Range("H7").Value = 5000
$wb.Worksheets.Item("BSM").Range("I7").Value = 5000
$wb.Worksheets.Item("BSM").Range("K7").Value = 5000
$wb.Worksheets.Item("BSM").Range("M7").Value = -4887

$wb.Worksheets.Item("BSM").Range("H60").Value = 69000
$wb.Worksheets.Item("BSM").Range("J60").Value = 69000
$wb.Worksheets.Item("BSM").Range("L60").Value = 69000
$wb.Worksheets.Item("BSM").Range("N60").Value = -70198

$wb.Worksheets.Item("BSM").Range("H86").Value = 3445
$wb.Worksheets.Item("BSM").Range("I86").Value = 1681.7693
$wb.Worksheets.Item("BSM").Range("J86").Value = 11085.667
$wb.Worksheets.Item("BSM").Range("K86").Value = 1681.7693
$wb.Worksheets.Item("BSM").Range("L86").Value = 11085.667
$wb.Worksheets.Item("BSM").Range("M86").Value = -558.7692999999999
$wb.Worksheets.Item("BSM").Range("N86").Value = -13331.667

$wb.Worksheets.Item("BSM").Range("H89").Value = 3445
$wb.Worksheets.Item("BSM").Range("I89").Value = 1681.7693
$wb.Worksheets.Item("BSM").Range("J89").Value = 11085.667
$wb.Worksheets.Item("BSM").Range("K89").Value = 8408.8465
$wb.Worksheets.Item("BSM").Range("L89").Value = 55428.335
$wb.Worksheets.Item("BSM").Range("M89").Value = -2792.8465
$wb.Worksheets.Item("BSM").Range("N89").Value = -66660.33499999999

$wb.Worksheets.Item("BSM").Range("H94").Value = 3511.8572
$wb.Worksheets.Item("BSM").Range("I94").Value = 491.73334
$wb.Worksheets.Item("BSM").Range("J94").Value = 6996.615
$wb.Worksheets.Item("BSM").Range("K94").Value = 491.73334
$wb.Worksheets.Item("BSM").Range("L94").Value = 6996.615
$wb.Worksheets.Item("BSM").Range("M94").Value = -40.73334
$wb.Worksheets.Item("BSM").Range("N94").Value = -7898.615

$wb.Worksheets.Item("BSM").Range("H100").Value = 29821.5
$wb.Worksheets.Item("BSM").Range("J100").Value = 29821.5
$wb.Worksheets.Item("BSM").Range("L100").Value = 29821.5
$wb.Worksheets.Item("BSM").Range("N100").Value = -31985.5

$wb.Worksheets.Item("CRP").Range("H2").Value = 0
$wb.Worksheets.Item("CRP").Range("J2").Value = 0
$wb.Worksheets.Item("CRP").Range("L2").Value = 0
$wb.Worksheets.Item("CRP").Range("N2").Value = $null

$wb.Worksheets.Item("CRP").Range("H3").Value = 18400
$wb.Worksheets.Item("CRP").Range("J3").Value = 1000
$wb.Worksheets.Item("CRP").Range("L3").Value = 1000
$wb.Worksheets.Item("CRP").Range("N3").Value = -1226

$wb.Worksheets.Item("CRP").Range("H7").Value = 184.73334
$wb.Worksheets.Item("CRP").Range("I7").Value = 184.73334
$wb.Worksheets.Item("CRP").Range("K7").Value = 184.73334
$wb.Worksheets.Item("CRP").Range("M7").Value = -71.73334

$wb.Worksheets.Item("CRP").Range("H22").Value = 909.6667
$wb.Worksheets.Item("CRP").Range("I22").Value = 813
$wb.Worksheets.Item("CRP").Range("J22").Value = 1248
$wb.Worksheets.Item("CRP").Range("K22").Value = 813
$wb.Worksheets.Item("CRP").Range("L22").Value = 1248
$wb.Worksheets.Item("CRP").Range("M22").Value = -463
$wb.Worksheets.Item("CRP").Range("N22").Value = -1948

$wb.Worksheets.Item("CRP").Range("H28").Value = 23912.375
$wb.Worksheets.Item("CRP").Range("J28").Value = 23912.375
$wb.Worksheets.Item("CRP").Range("L28").Value = 23912.375
$wb.Worksheets.Item("CRP").Range("N28").Value = -24402.375

$wb.Worksheets.Item("CRP").Range("H31").Value = 5332.2104
$wb.Worksheets.Item("CRP").Range("I31").Value = 1211.2858
$wb.Worksheets.Item("CRP").Range("J31").Value = 7736.0835
$wb.Worksheets.Item("CRP").Range("K31").Value = 1211.2858
$wb.Worksheets.Item("CRP").Range("L31").Value = 7736.0835
$wb.Worksheets.Item("CRP").Range("M31").Value = -916.2858000000001
$wb.Worksheets.Item("CRP").Range("N31").Value = -8326.083500000001

$wb.Worksheets.Item("CRP").Range("H34").Value = 5332.2104
$wb.Worksheets.Item("CRP").Range("I34").Value = 1211.2858
$wb.Worksheets.Item("CRP").Range("J34").Value = 7736.0835
$wb.Worksheets.Item("CRP").Range("K34").Value = 1211.2858
$wb.Worksheets.Item("CRP").Range("L34").Value = 7736.0835
$wb.Worksheets.Item("CRP").Range("M34").Value = -1009.2858
$wb.Worksheets.Item("CRP").Range("N34").Value = -8140.0835

$wb.Worksheets.Item("CRP").Range("H58").Value = 17786.777
$wb.Worksheets.Item("CRP").Range("I58").Value = 14402.2
$wb.Worksheets.Item("CRP").Range("K58").Value = 14402.2
$wb.Worksheets.Item("CRP").Range("M58").Value = -14199.2

$wb.Worksheets.Item("CRP").Range("H62").Value = 12310.8
$wb.Worksheets.Item("CRP").Range("I62").Value = 4202.5
$wb.Worksheets.Item("CRP").Range("K62").Value = 4202.5
$wb.Worksheets.Item("CRP").Range("M62").Value = -3578.5

$wb.Worksheets.Item("CRP").Range("H65").Value = 12310.8
$wb.Worksheets.Item("CRP").Range("I65").Value = 4202.5
$wb.Worksheets.Item("CRP").Range("K65").Value = 21012.5
$wb.Worksheets.Item("CRP").Range("M65").Value = -17892.5

$wb.Worksheets.Item("CRP").Range("H107").Value = 961.84
$wb.Worksheets.Item("CRP").Range("I107").Value = 702.25
$wb.Worksheets.Item("CRP").Range("K107").Value = 702.25
$wb.Worksheets.Item("CRP").Range("M107").Value = 1217.75

$wb.Worksheets.Item("CRP").Range("H132").Value = 6605.72
$wb.Worksheets.Item("CRP").Range("I132").Value = 6211.6665
$wb.Worksheets.Item("CRP").Range("J132").Value = 7619
$wb.Worksheets.Item("CRP").Range("K132").Value = 18634.9995
$wb.Worksheets.Item("CRP").Range("L132").Value = 22857
$wb.Worksheets.Item("CRP").Range("M132").Value = -16104.9995
$wb.Worksheets.Item("CRP").Range("N132").Value = -27917

$wb.Worksheets.Item("CRP").Range("H134").Value = 602832.9
$wb.Worksheets.Item("CRP").Range("I134").Value = 1493566.1
$wb.Worksheets.Item("CRP").Range("K134").Value = 4480698.300000001
$wb.Worksheets.Item("CRP").Range("M134").Value = -4478163.300000001

$wb.Worksheets.Item("CRP").Range("H136").Value = 17786.777
$wb.Worksheets.Item("CRP").Range("I136").Value = 14402.2
$wb.Worksheets.Item("CRP").Range("K136").Value = 43206.60000000001
$wb.Worksheets.Item("CRP").Range("M136").Value = -40656.60000000001

$wb.Worksheets.Item("CUL").Range("H2").Value = 139.8
$wb.Worksheets.Item("CUL").Range("J2").Value = 150
$wb.Worksheets.Item("CUL").Range("L2").Value = 900
$wb.Worksheets.Item("CUL").Range("N2").Value = -1126

$wb.Worksheets.Item("CUL").Range("H12").Value = 1176851.2
$wb.Worksheets.Item("CUL").Range("I12").Value = 2500043
$wb.Worksheets.Item("CUL").Range("K12").Value = 7500129
$wb.Worksheets.Item("CUL").Range("M12").Value = -7499956

$wb.Worksheets.Item("CUL").Range("H23").Value = 56.8
$wb.Worksheets.Item("CUL").Range("I23").Value = 26
$wb.Worksheets.Item("CUL").Range("J23").Value = 77.333336
$wb.Worksheets.Item("CUL").Range("K23").Value = 78
$wb.Worksheets.Item("CUL").Range("L23").Value = 232.000008
$wb.Worksheets.Item("CUL").Range("M23").Value = 157
$wb.Worksheets.Item("CUL").Range("N23").Value = -702.000008

$wb.Worksheets.Item("CUL").Range("H38").Value = 37.416668
$wb.Worksheets.Item("CUL").Range("J38").Value = 64.5
$wb.Worksheets.Item("CUL").Range("L38").Value = 193.5
$wb.Worksheets.Item("CUL").Range("N38").Value = -887.5

$wb.Worksheets.Item("CUL").Range("H122").Value = 138512.38
$wb.Worksheets.Item("CUL").Range("I122").Value = 598.75
$wb.Worksheets.Item("CUL").Range("J122").Value = 155229.19
$wb.Worksheets.Item("CUL").Range("K122").Value = 5388.75
$wb.Worksheets.Item("CUL").Range("L122").Value = 1397062.71
$wb.Worksheets.Item("CUL").Range("M122").Value = -2938.75
$wb.Worksheets.Item("CUL").Range("N122").Value = -1401962.71

$wb.Worksheets.Item("CUL").Range("H131").Value = 41272308
$wb.Worksheets.Item("CUL").Range("I131").Value = 106667260
$wb.Worksheets.Item("CUL").Range("J131").Value = 20836382
$wb.Worksheets.Item("CUL").Range("K131").Value = 320001780
$wb.Worksheets.Item("CUL").Range("L131").Value = 62509146
$wb.Worksheets.Item("CUL").Range("M131").Value = -319996740
$wb.Worksheets.Item("CUL").Range("N131").Value = -62519226

$wb.Worksheets.Item("GSM").Range("H40").Value = 20000
$wb.Worksheets.Item("GSM").Range("I40").Value = 20000
$wb.Worksheets.Item("GSM").Range("K40").Value = 20000
$wb.Worksheets.Item("GSM").Range("M40").Value = -19849

$wb.Worksheets.Item("LTW").Range("H20").Value = 26800
$wb.Worksheets.Item("LTW").Range("I20").Value = 25600
$wb.Worksheets.Item("LTW").Range("K20").Value = 25600
$wb.Worksheets.Item("LTW").Range("M20").Value = -25374

$wb.Worksheets.Item("LTW").Range("H22").Value = 1083.1666
$wb.Worksheets.Item("LTW").Range("I22").Value = 1000
$wb.Worksheets.Item("LTW").Range("J22").Value = 1099.8
$wb.Worksheets.Item("LTW").Range("K22").Value = 1000
$wb.Worksheets.Item("LTW").Range("L22").Value = 1099.8
$wb.Worksheets.Item("LTW").Range("M22").Value = -705
$wb.Worksheets.Item("LTW").Range("N22").Value = -1689.8

$wb.Worksheets.Item("LTW").Range("H27").Value = 1083.1666
$wb.Worksheets.Item("LTW").Range("I27").Value = 1000
$wb.Worksheets.Item("LTW").Range("J27").Value = 1099.8
$wb.Worksheets.Item("LTW").Range("K27").Value = 1000
$wb.Worksheets.Item("LTW").Range("L27").Value = 1099.8
$wb.Worksheets.Item("LTW").Range("M27").Value = -893
$wb.Worksheets.Item("LTW").Range("N27").Value = -1313.8

$wb.Worksheets.Item("LTW").Range("H61").Value = 14260.223
$wb.Worksheets.Item("LTW").Range("J61").Value = 21161.25
$wb.Worksheets.Item("LTW").Range("L61").Value = 21161.25
$wb.Worksheets.Item("LTW").Range("N61").Value = -21565.25

$wb.Worksheets.Item("LTW").Range("H104").Value = 31946.727
$wb.Worksheets.Item("LTW").Range("J104").Value = 31946.727
$wb.Worksheets.Item("LTW").Range("L104").Value = 31946.727
$wb.Worksheets.Item("LTW").Range("N104").Value = -38934.727

$wb.Worksheets.Item("LTW").Range("H113").Value = 14260.223
$wb.Worksheets.Item("LTW").Range("J113").Value = 21161.25
$wb.Worksheets.Item("LTW").Range("L113").Value = 21161.25
$wb.Worksheets.Item("LTW").Range("N113").Value = -25501.25

$wb.Worksheets.Item("LTW").Range("H122").Value = 4236.8125
$wb.Worksheets.Item("LTW").Range("I122").Value = 3483.7693
$wb.Worksheets.Item("LTW").Range("J122").Value = 7500
$wb.Worksheets.Item("LTW").Range("K122").Value = 10451.3079
$wb.Worksheets.Item("LTW").Range("L122").Value = 22500
$wb.Worksheets.Item("LTW").Range("M122").Value = -8001.3079
$wb.Worksheets.Item("LTW").Range("N122").Value = -27400

$wb.Worksheets.Item("WVR").Range("H122").Value = 3150.1482
$wb.Worksheets.Item("WVR").Range("I122").Value = 3150.1482
$wb.Worksheets.Item("WVR").Range("K122").Value = 9450.444600000001
$wb.Worksheets.Item("WVR").Range("M122").Value = -7000.444600000001

$wb.Worksheets.Item("WVR").Range("H126").Value = 4371.2
$wb.Worksheets.Item("WVR").Range("I126").Value = 2986.5
$wb.Worksheets.Item("WVR").Range("K126").Value = 8959.5
$wb.Worksheets.Item("WVR").Range("M126").Value = -6489.5
